$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.537.78'
$ws.Range('E2').Value = '  +2.88%  '
$ws.Range('D3').Value = '2.313.56'
$ws.Range('E3').Value = '  +1.80%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').Value = '''311.31'
$ws.Range('E5').Value = '  +1.67%  '
$ws.Range('D6').Value = '''102.32'
$ws.Range('E6').Value = '  +4.65%  '
$ws.Range('D7').Value = '''0.538'
$ws.Range('E7').Value = '  +1.84%  '
$ws.Range('E8').Value = '  -0.03%  '
$ws.Range('E9').Value = '  +7.54%  '
$ws.Range('D10').Value = '''35.86'
$ws.Range('E10').Value = '  +1.53%  '
$ws.Range('D11').Value = '''0.0818'
$ws.Range('E11').Value = '  +3.52%  '
$ws.Range('D12').Value = '''0.113'
$ws.Range('E12').Value = '  -0.32%  '
$ws.Range('E13').Value = '  +1.17%  '
$ws.Range('D14').Value = '2.670.52'
$ws.Range('E14').Value = '  +1.71%  '
$ws.Range('D15').Value = '''15.01'
$ws.Range('E15').Value = '  +1.85%  '
$ws.Range('D16').Value = '2.311.96'
$ws.Range('E16').Value = '  +2.49%  '
$ws.Range('E17').Value = '  +1.97%  '
$ws.Range('D18').Value = '43.425.75'
$ws.Range('E18').Value = '  +2.96%  '
$ws.Range('D19').Value = '''12.44'
$ws.Range('E19').Value = '  -0.30%  '
$ws.Range('D20').Value = '0.0₃0929'
$ws.Range('E20').Value = '  +2.46%  '
$ws.Range('D21').Value = '''6.18'
$ws.Range('E21').Value = '  +2.51%  '
$ws.Range('D22').Value = '''68.21'
$ws.Range('E22').Value = '  +0.22%  '
$ws.Range('D23').Value = '''242.11'
$ws.Range('E24').Value = '  +2.84%  '
$ws.Range('D25').Value = '''2.63'
$ws.Range('E25').Value = '  +2.07%  '
$ws.Range('E26').Value = '  -0.01%  '
$ws.Range('D27').Value = '''3.99'
$ws.Range('E27').Value = '  -1.55%  '
$ws.Range('D28').Value = '''24.78'
$ws.Range('E28').Value = '  +4.81%  '
$ws.Range('E29').Value = '  +8.05%  '
$ws.Range('D30').Value = '''36.81'
$ws.Range('E30').Value = '  -2.77%  '
$ws.Range('D31').Value = '''9.66'
$ws.Range('E31').Value = '  +1.34%  '
$ws.Range('D32').Value = '''167.06'
$ws.Range('E32').Value = '  +3.52%  '
$ws.Range('D33').Value = '''5.30'
$ws.Range('E33').Value = '  +1.02%  '
$ws.Range('E34').Value = '  +0.09%  '
$ws.Range('B35').Value = 'Hedera'
$ws.Range('C35').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D35').Value = '''0.0746'
$ws.Range('E35').Value = '  +0.99%  '
$ws.Range('B36').Value = 'LidoDAOToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D36').Value = '''3.10'
$ws.Range('E36').Value = '  -1.83%  '
$ws.Range('E37').Value = '  +5.68%  '
$ws.Range('D38').Value = '''17.63'
$ws.Range('E38').Value = '  -0.45%  '
$ws.Range('D39').Value = '''0.107'
$ws.Range('E39').Value = '  +1.58%  '
$ws.Range('D40').Value = '''1.88'
$ws.Range('E40').Value = '  +2.85%  '
$ws.Range('E41').Value = '  +1.71%  '
$ws.Range('D42').Value = '''4.32'
$ws.Range('E42').Value = '  +5.67%  '
$ws.Range('E43').Value = '  -1.00%  '
$ws.Range('D44').Value = '''19.37'
$ws.Range('E44').Value = '  +2.43%  '
$ws.Range('E45').Value = '  +2.94%  '
$ws.Range('D46').Value = '1.971.06'
$ws.Range('E46').Value = '  +1.21%  '
$ws.Range('D47').Value = '''2.99'
$ws.Range('E47').Value = '  +2.58%  '
$ws.Range('D48').Value = '''10.02'
$ws.Range('E48').Value = '  +0.97%  '
$ws.Range('B49').Value = 'HuobiToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D49').Value = '''2.95'
$ws.Range('E49').Value = '  +6.00%  '
$ws.Range('B50').Value = 'MultiversX'
$ws.Range('C50').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D50').Value = '''55.62'
$ws.Range('E50').Value = '  +3.56%  '
$ws.Range('D51').Value = '''1.58'
$ws.Range('E51').Value = '  +7.14%  '
